$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map old section labels to the new zero-padded labels.
$map = @{
    "CSE-1" = "CSE-01"
    "CSE-2" = "CSE-02"
    "CSE-3" = "CSE-03"
    "CSE-4" = "CSE-04"
    "CSE-5" = "CSE-05"
    "CSE-6" = "CSE-06"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}

# Update the saved view/selection state of the sheet.
$ws.Range("P144").Select()
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 133
$activeWindow.ScrollColumn = 1
